# Add a new "TB_CALL_HISTORY" worksheet between TB_PROJECTS and TB_CODE,
# describing a call/consultation-history table, and nudge a couple of
# cell selections on the other sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the new sheet by duplicating TB_CODE (same look & feel: same
#    header style, same border/row style for data cells) and placing the
#    duplicate immediately before TB_CODE - exactly where the new table
#    belongs in the sheet order.
# ---------------------------------------------------------------------
$codeSheet = $wb.Worksheets.Item("TB_CODE")
$codeSheet.Copy($codeSheet, $null)
$newSheet = $wb.Worksheets.Item("TB_CODE (2)")
$newSheet.Name = "TB_CALL_HISTORY"

# ---------------------------------------------------------------------
# 2. Fill in the table-specification content for TB_CALL_HISTORY.
#    Columns: Column Name | Column Id | Data Type | Constraints | Description
#    Filled Column Id + Data Type first (row by row), then Column Name,
#    then Constraints - same authoring order as the source file, so new
#    shared-string entries land in the same sequence.
# ---------------------------------------------------------------------

# Rows 6-10 did not exist on the copied TB_CODE sheet (it only had 5
# rows) - copy the row-5 formatting down first so the new rows pick up
# the same bordered, non-wrapped data-row style, before filling values.
$newSheet.Range("A5:E5").Copy()
$newSheet.Range("A6:E10").PasteSpecial(-4122)
$newSheet.Range("A1").Select()

$newSheet.Range("B2").Value = "REG_DATE"
$newSheet.Range("C2").Value = "VARCHAR2(8 BYTE)"

$newSheet.Range("B3").Value = "EMPLOYEE_NAME"
$newSheet.Range("C3").Value = "VARCHAR2(100 BYTE)"

$newSheet.Range("B4").Value = "HISTORY_ID"
$newSheet.Range("C4").Value = "NUMBER"

$newSheet.Range("B5").Value = "NOTE"
$newSheet.Range("C5").Value = "VARCHAR2(500 BYTE)"

$newSheet.Range("B6").Value = "ACTION"
$newSheet.Range("C6").Value = "VARCHAR2(500 BYTE)"

$newSheet.Range("B7").Value = "REGDATE"
$newSheet.Range("C7").Value = "TIMESTAMP(6)"

$newSheet.Range("B8").Value = "REGID"
$newSheet.Range("C8").Value = "VARCHAR2(100 BYTE)"

$newSheet.Range("B9").Value = "MODDATE"
$newSheet.Range("C9").Value = "TIMESTAMP(6)"

$newSheet.Range("B10").Value = "MODID"
$newSheet.Range("C10").Value = "VARCHAR2(100 BYTE)"

$newSheet.Range("A2").Value = "등록일자"
$newSheet.Range("A3").Value = "개발자명"
$newSheet.Range("A4").Value = "상담이력ID"
$newSheet.Range("A5").Value = "상담내용"
$newSheet.Range("A6").Value = "조치내용"
$newSheet.Range("A7").Value = "등록일시"
$newSheet.Range("A8").Value = "등록자"
$newSheet.Range("A9").Value = "수정일시"
$newSheet.Range("A10").Value = "수정자"

$newSheet.Range("D2").Value = "NOT NULL"
$newSheet.Range("D3").Value = "NOT NULL"
$newSheet.Range("D4").Value = "NOT NULL"
$newSheet.Range("D5").Value = ""

# ---------------------------------------------------------------------
# 3. Data rows on TB_CODE used the wrapping "NOT NULL" style for every
#    row; on the new sheet only columns A-C (always) and D (for blank
#    constraint cells) drop the word-wrap so they match the rest of the
#    workbook's plain data-row look.
# ---------------------------------------------------------------------
$newSheet.Range("A2:C10").WrapText = $false
$newSheet.Range("E2:E10").WrapText = $false
$newSheet.Range("D5:D10").WrapText = $false

# Header row is taller on this sheet (two-line wrapped header).
$newSheet.Rows.Item(1).RowHeight = 34.8

$newSheet.Range("A1:E10").Select()
$newSheet.Range("A1").Select()

# ---------------------------------------------------------------------
# 4. Cursor/selection nudges recorded on the other sheets.
# ---------------------------------------------------------------------
$empSheet = $wb.Worksheets.Item("TB_EMPLOYEES")
$empSheet.Activate()
$empSheet.Range("E17").Select()

$projSheet = $wb.Worksheets.Item("TB_PROJECTS")
$projSheet.Activate()
$projSheet.Range("D14").Select()

$newSheet.Activate()
$newSheet.Range("C14").Select()

# Re-fetch TB_CODE by name: the original $codeSheet reference was taken
# before the sheet list was reshuffled by Copy(), so look it up fresh.
$codeSheet2 = $wb.Worksheets.Item("TB_CODE")
$codeSheet2.Activate()
$codeSheet2.PageSetup.PaperSize = 9
$codeSheet2.PageSetup.Orientation = 1

# TB_PROJECTS stays the active tab, matching the saved workbook state.
$projSheet.Activate()
